$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.25
$ws.Range("C2").Value = 0.5
$ws.Range("J2").Value = 0.125
$ws.Range("S2").Value = 0.125
$ws.Range("P3").Value = 0.5
$ws.Range("S3").Value = 0.5
$ws.Range("P4").Value = 0.3333333333333333
$ws.Range("S4").Value = 0.6666666666666666
$ws.Range("F6").Value = 0.1153846153846154
$ws.Range("J6").Value = 0.3461538461538461
$ws.Range("O6").Value = 0.07692307692307693
$ws.Range("S6").Value = 0.4615384615384616
$ws.Range("B7").Value = 0.04761904761904762
$ws.Range("F7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.09523809523809523
$ws.Range("Q7").Value = 0.1904761904761905
$ws.Range("R7").Value = 0.04761904761904762
$ws.Range("S7").Value = 0.5714285714285714
$ws.Range("B8").Value = 0.02777777777777778
$ws.Range("D8").Value = 0.02777777777777778
$ws.Range("F8").Value = 0.02777777777777778
$ws.Range("J8").Value = 0.1111111111111111
$ws.Range("Q8").Value = 0.1111111111111111
$ws.Range("R8").Value = 0.05555555555555555
$ws.Range("S8").Value = 0.6388888888888888
$ws.Range("F9").Value = 0.09090909090909091
$ws.Range("J9").Value = 0.09090909090909091
$ws.Range("Q9").Value = 0.1818181818181818
$ws.Range("R9").Value = 0.2727272727272727
$ws.Range("S9").Value = 0.3636363636363636
$ws.Range("B10").Value = 0.04651162790697674
$ws.Range("D10").Value = 0.02325581395348837
$ws.Range("F10").Value = 0.06976744186046512
$ws.Range("J10").Value = 0.1046511627906977
$ws.Range("O10").Value = 0.01162790697674419
$ws.Range("Q10").Value = 0.2325581395348837
$ws.Range("R10").Value = 0.1279069767441861
$ws.Range("S10").Value = 0.3837209302325582
$ws.Range("G11").Value = 0.125
$ws.Range("K11").Value = 0.125
$ws.Range("L11").Value = 0.75
$ws.Range("G12").Value = 0.9411764705882353
$ws.Range("J12").Value = 0.05882352941176471
$ws.Range("J13").Value = 0.1666666666666667
$ws.Range("S13").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.0625
$ws.Range("H15").Value = 0.125
$ws.Range("J15").Value = 0.5625
$ws.Range("K15").Value = 0.0625
$ws.Range("S15").Value = 0.1875
$ws.Range("H16").Value = 0.3333333333333333
$ws.Range("J16").Value = 0.3333333333333333
$ws.Range("K16").Value = 0.3333333333333333
$ws.Range("F17").Value = 0.1
$ws.Range("H17").Value = 0.1
$ws.Range("J17").Value = 0.3
$ws.Range("K17").Value = 0.1333333333333333
$ws.Range("O17").Value = 0.03333333333333333
$ws.Range("H18").Value = 0.1176470588235294
$ws.Range("I18").Value = 0.05882352941176471
$ws.Range("J18").Value = 0.2941176470588235
$ws.Range("K18").Value = 0.1176470588235294
$ws.Range("O18").Value = 0.2352941176470588
$ws.Range("S18").Value = 0.1764705882352941
$ws.Range("F19").Value = 0.06086956521739131
$ws.Range("H19").Value = 0.2434782608695652
$ws.Range("I19").Value = 0.04347826086956522
$ws.Range("J19").Value = 0.3130434782608696
$ws.Range("K19").Value = 0.1043478260869565
$ws.Range("M19").Value = 0.03478260869565217
$ws.Range("O19").Value = 0.06086956521739131
$ws.Range("S19").Value = 0.1391304347826087
